$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Brand, Family) after the "Product" column (A),
# shifting the old B..F columns to D..H.
$ws.Columns("B:C").Insert()

# Header row
$ws.Range("B1").Value = "Brand"
$ws.Range("C1").Value = "Family"

# Product 1 & Product 2 rows (2-16) -> Brand W / Family A
$ws.Range("B2:B16").Value = "Brand W"
$ws.Range("C2:C16").Value = "Family A"

# Product 3, 4, 5 rows (17-39) -> Brand X / Family A
$ws.Range("B17:B39").Value = "Brand X"
$ws.Range("C17:C39").Value = "Family A"

# View settings: zoom to 150% and update the active selection/window.
$win = $excel.ActiveWindow
$win.Zoom = 150
$ws.Range("B10").Select()
